$wb = $excel.ActiveWorkbook

# Scheduled-runner price/profit refresh for the Leve profit tables.
# Updates currentAveragePrice(NQ/HQ) (H:J), LevePrice(NQ/HQ) (K:L) and the
# resulting LeveProfit(NQ/HQ) (M:N) on each job sheet to the latest market
# data. Some leves that are no longer profitable/tracked have their
# profit columns cleared instead of recomputed.

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H116").Value = 8549441
$ws.Range("I116").Value = 19232968
$ws.Range("K116").Value = 19232968
$ws.Range("M116").Value = -19229526
$ws.Range("H132").Value = 3224.0889
$ws.Range("I132").Value = 3192.3447
$ws.Range("J132").Value = 3281.625
$ws.Range("K132").Value = 9577.034100000001
$ws.Range("L132").Value = 9844.875
$ws.Range("M132").Value = -7047.034100000001
$ws.Range("N132").Value = -14904.875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1659.3334
$ws.Range("I61").Value = 1158.2084
$ws.Range("J61").Value = 3663.8333
$ws.Range("K61").Value = 1158.2084
$ws.Range("L61").Value = 3663.8333
$ws.Range("M61").Value = -946.2084
$ws.Range("N61").Value = -4087.8333
$ws.Range("H74").Value = 3326.9285
$ws.Range("I74").Value = 3616.7778
$ws.Range("J74").Value = 2805.2
$ws.Range("K74").Value = 3616.7778
$ws.Range("L74").Value = 2805.2
$ws.Range("M74").Value = -2742.7778
$ws.Range("N74").Value = -4553.2
$ws.Range("H77").Value = 3326.9285
$ws.Range("I77").Value = 3616.7778
$ws.Range("J77").Value = 2805.2
$ws.Range("K77").Value = 18083.889
$ws.Range("L77").Value = 14026
$ws.Range("M77").Value = -13715.889
$ws.Range("N77").Value = -22762
$ws.Range("H102").Value = 1871.25
$ws.Range("I102").Value = 1624.8334
$ws.Range("J102").Value = 2610.5
$ws.Range("K102").Value = 1624.8334
$ws.Range("L102").Value = 2610.5
$ws.Range("M102").Value = -2.833399999999983
$ws.Range("N102").Value = -5854.5
$ws.Range("H122").Value = 902.62067
$ws.Range("I122").Value = 821.7037
$ws.Range("K122").Value = 2465.1111
$ws.Range("M122").Value = -15.11110000000008
$ws.Range("H132").Value = 1483
$ws.Range("I132").Value = 1337.8948
$ws.Range("K132").Value = 4013.6844
$ws.Range("M132").Value = -1483.6844
$ws.Range("H136").Value = 1659.3334
$ws.Range("I136").Value = 1158.2084
$ws.Range("J136").Value = 3663.8333
$ws.Range("K136").Value = 3474.6252
$ws.Range("L136").Value = 10991.4999
$ws.Range("M136").Value = -924.6251999999999
$ws.Range("N136").Value = -16091.4999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H99").Value = 9380
$ws.Range("I99").Value = 20677.6
$ws.Range("J99").Value = 3103.5557
$ws.Range("K99").Value = 20677.6
$ws.Range("L99").Value = 3103.5557
$ws.Range("M99").Value = -19179.6
$ws.Range("N99").Value = -6099.5557
$ws.Range("H105").Value = 2831.39
$ws.Range("I105").Value = 1993.3334
$ws.Range("J105").Value = 2979.2825
$ws.Range("K105").Value = 1993.3334
$ws.Range("L105").Value = 2979.2825
$ws.Range("M105").Value = -246.3334
$ws.Range("N105").Value = -6473.282499999999
$ws.Range("H122").Value = 35380
$ws.Range("J122").Value = 35380
$ws.Range("L122").Value = 35380
$ws.Range("N122").Value = -45180
$ws.Range("H134").Value = 1742.9423
$ws.Range("I134").Value = 1539.85
$ws.Range("K134").Value = 4619.549999999999
$ws.Range("M134").Value = -2084.549999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5532.385
$ws.Range("I86").Value = 7485.5
$ws.Range("J86").Value = 3858.2856
$ws.Range("K86").Value = 7485.5
$ws.Range("L86").Value = 3858.2856
$ws.Range("M86").Value = -6362.5
$ws.Range("N86").Value = -6104.2856
$ws.Range("H89").Value = 5532.385
$ws.Range("I89").Value = 7485.5
$ws.Range("J89").Value = 3858.2856
$ws.Range("K89").Value = 37427.5
$ws.Range("L89").Value = 19291.428
$ws.Range("M89").Value = -31811.5
$ws.Range("N89").Value = -30523.428
$ws.Range("H105").Value = 664.2
$ws.Range("I105").Value = 736.6667
$ws.Range("J105").Value = 555.5
$ws.Range("K105").Value = 736.6667
$ws.Range("L105").Value = 555.5
$ws.Range("M105").Value = 1010.3333
$ws.Range("N105").Value = -4049.5
$ws.Range("H132").Value = 2095.32
$ws.Range("I132").Value = 1852.5625
$ws.Range("J132").Value = 2526.889
$ws.Range("K132").Value = 5557.6875
$ws.Range("L132").Value = 7580.667
$ws.Range("M132").Value = -3027.6875
$ws.Range("N132").Value = -12640.667
$ws.Range("H134").Value = 2010.6046
$ws.Range("I134").Value = 1284.7878
$ws.Range("J134").Value = 4405.8
$ws.Range("K134").Value = 3854.3634
$ws.Range("L134").Value = 13217.4
$ws.Range("M134").Value = -1319.3634
$ws.Range("N134").Value = -18287.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 2994.4546
$ws.Range("I99").Value = 2994.4546
$ws.Range("K99").Value = 2994.4546
$ws.Range("M99").Value = -748.4546
$ws.Range("H122").Value = 909.1875
$ws.Range("I122").Value = 922.4545000000001
$ws.Range("J122").Value = 880
$ws.Range("K122").Value = 2767.3635
$ws.Range("L122").Value = 2640
$ws.Range("M122").Value = -317.3635000000004
$ws.Range("N122").Value = -7540
$ws.Range("H132").Value = 3712.5
$ws.Range("I132").Value = 3806.8157
$ws.Range("J132").Value = 3354.1
$ws.Range("K132").Value = 11420.4471
$ws.Range("L132").Value = 10062.3
$ws.Range("M132").Value = -8890.447100000001
$ws.Range("N132").Value = -15122.3

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 72860990
$ws.Range("I100").Value = 2503495
$ws.Range("J100").Value = 166670990
$ws.Range("K100").Value = 2503495
$ws.Range("L100").Value = 166670990
$ws.Range("M100").Value = -166672072

Write-Output "Applied scheduled runner updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets."
